$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Stash the three existing row styles (none / green / orange) into
#    throwaway helper cells, far outside the used range, before any
#    other edits happen. They are used below as stable PasteSpecial
#    format donors (so the tables fills reuse the workbooks
#    existing style slots instead of creating duplicate styles), and
#    are cleared again at the very end.
# ------------------------------------------------------------------
$ws.Range("A8:B8").Copy() | Out-Null        # donor: no fill
$ws.Range("Z1:AA1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:B2").Copy() | Out-Null        # donor: green fill
$ws.Range("Z2:AA2").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:B5").Copy() | Out-Null        # donor: orange fill
$ws.Range("Z3:AA3").PasteSpecial(-4122) | Out-Null
$noneDonor = "Z1:AA1"
$greenDonor = "Z2:AA2"
$orangeDonor = "Z3:AA3"

# ------------------------------------------------------------------
# 2) Grow the table to 3 columns x 22 rows.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C22")) | Out-Null

# ------------------------------------------------------------------
# 3) Re-apply the correct fill to every data row (A:C) BEFORE rewriting
#    values, so each row ends up with the right style index.
# ------------------------------------------------------------------
$ws.Range($noneDonor).Copy() | Out-Null
$ws.Range('A2:C2').PasteSpecial(-4122) | Out-Null  # none
$ws.Range($noneDonor).Copy() | Out-Null
$ws.Range('A13:C13').PasteSpecial(-4122) | Out-Null  # none
$ws.Range($noneDonor).Copy() | Out-Null
$ws.Range('A15:C15').PasteSpecial(-4122) | Out-Null  # none
$ws.Range($noneDonor).Copy() | Out-Null
$ws.Range('A17:C21').PasteSpecial(-4122) | Out-Null  # none
$ws.Range($greenDonor).Copy() | Out-Null
$ws.Range('A3:C6').PasteSpecial(-4122) | Out-Null  # green
$ws.Range($greenDonor).Copy() | Out-Null
$ws.Range('A9:C12').PasteSpecial(-4122) | Out-Null  # green
$ws.Range($greenDonor).Copy() | Out-Null
$ws.Range('A14:C14').PasteSpecial(-4122) | Out-Null  # green
$ws.Range($greenDonor).Copy() | Out-Null
$ws.Range('A16:C16').PasteSpecial(-4122) | Out-Null  # green
$ws.Range($greenDonor).Copy() | Out-Null
$ws.Range('A22:C22').PasteSpecial(-4122) | Out-Null  # green
$ws.Range($orangeDonor).Copy() | Out-Null
$ws.Range('A7:C8').PasteSpecial(-4122) | Out-Null  # orange
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Write the final cell values (A: task, B: priority, C: comment).
# ------------------------------------------------------------------
$ws.Range('A1').Value2 = 'Что сделать'
$ws.Range('B1').Value2 = 'Приоритет'
$ws.Range('C1').Value2 = 'Комментарий'
$ws.Range('A2').Value2 = 'Увеличить количество разрешенных запросов до 100 в течении 1 мин и блокировкой на 15 сек'
$ws.Range('B2').Value2 = 15
$ws.Range('A3').Value2 = 'Зайти в 1 курс, перейти в другой курс, отображается урок предыдущего курса'
$ws.Range('B3').Value2 = 10
$ws.Range('A4').Value2 = 'Реализовать filter для courses'
$ws.Range('B4').Value2 = 10
$ws.Range('A5').Value2 = 'Реализовать добавление в избранное'
$ws.Range('B5').Value2 = 10
$ws.Range('A6').Value2 = 'Сделать адаптив'
$ws.Range('B6').Value2 = 10
$ws.Range('A7').Value2 = 'сделать перелистывание страниц courses. Отображать по 10 курсов на странице'
$ws.Range('B7').Value2 = 9
$ws.Range('C7').Value2 = 'Нет функционала api'
$ws.Range('A8').Value2 = 'Выдавать сертификат только после успешного выполнения последнего урока. Изменить api. Сертификат в course/id И появляется после выполнения последнего теста'
$ws.Range('B8').Value2 = 9
$ws.Range('C8').Value2 = 'Нет функционала api'
$ws.Range('A9').Value2 = 'При клике в меню на модуль 2, а потом на модуль 1, выскакивает ошибка'
$ws.Range('B9').Value2 = 8
$ws.Range('A10').Value2 = 'Убрать все ошибки и warning'
$ws.Range('B10').Value2 = 8
$ws.Range('A11').Value2 = 'Протестировать каждый запрос к серверу на вечный цикл'
$ws.Range('B11').Value2 = 7
$ws.Range('A12').Value2 = 'В конце убрать все debugger, console.log'
$ws.Range('B12').Value2 = 6
$ws.Range('A13').Value2 = 'Убрать дублирование кода'
$ws.Range('B13').Value2 = 5
$ws.Range('A14').Value2 = 'Везде вместо "Подождите, идет загрузка", поставить прелоадер. Но чтобы не было несколько прелоадеров друг над другом'
$ws.Range('B14').Value2 = 5
$ws.Range('A15').Value2 = 'Пофиксить баг - при переходе из модуля 2 в модуль 1 по кнопке предыдущий урок, перекидывает на 1 урок 1 модуля'
$ws.Range('B15').Value2 = 4
$ws.Range('A16').Value2 = 'Удалить history и qs библиотеки'
$ws.Range('B16').Value2 = 4
$ws.Range('A17').Value2 = 'Переименовать файлы логично'
$ws.Range('B17').Value2 = 3
$ws.Range('A18').Value2 = 'Оптимизировать css'
$ws.Range('B18').Value2 = 3
$ws.Range('A19').Value2 = 'Удалить неиспользуемые props'
$ws.Range('B19').Value2 = 3
$ws.Range('A20').Value2 = 'Исправить желтые подчеркивания'
$ws.Range('B20').Value2 = 3
$ws.Range('A21').Value2 = 'Убрать выделение мышью там, где могут несколько раз кликнуть'
$ws.Range('B21').Value2 = 2
$ws.Range('A22').Value2 = 'Удалить неиспользуемые imports'
$ws.Range('B22').Value2 = 2

# ------------------------------------------------------------------
# 5) Clean up helper cells, size the new column, set the selection.
# ------------------------------------------------------------------
$ws.Range("Z1:AA3").Clear() | Out-Null
$ws.Columns("C").ColumnWidth = 19.33
$ws.Range("A13").Select() | Out-Null

Write-Output ("table range: " + $lo.Range.Address())
